$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44343
$ws.Range("H2").Value = 'Americana (o)'
$ws.Range("I2").Value = 'Primera'
$ws.Range("J2").Value = 20
$ws.Range("K2").Value = 36000
$ws.Range("L2").Value = 36000
$ws.Range("M2").Value = 36000
$ws.Range("N2").Value = '$/caja 25 kilos'
$ws.Range("P2").Value = 1440
$ws.Range("Q2").Value = 25
$ws.Range("D3").Value = 44221
$ws.Range("H3").Value = 'Americana (o)'
$ws.Range("I3").Value = 'Primera'
$ws.Range("J3").Value = 22
$ws.Range("K3").Value = 24000
$ws.Range("L3").Value = 25000
$ws.Range("M3").Value = 24545
$ws.Range("N3").Value = '$/caja 25 kilos'
$ws.Range("P3").Value = 982
$ws.Range("Q3").Value = 25
$ws.Range("D4").Value = 44581
$ws.Range("H4").Value = 'Americana (o)'
$ws.Range("I4").Value = 'Segunda'
$ws.Range("J4").Value = 30
$ws.Range("K4").Value = 17000
$ws.Range("L4").Value = 17000
$ws.Range("M4").Value = 17000
$ws.Range("N4").Value = '$/caja 25 kilos'
$ws.Range("P4").Value = 680
$ws.Range("Q4").Value = 25
$ws.Range("D5").Value = 44319
$ws.Range("H5").Value = 'Americana (o)'
$ws.Range("I5").Value = 'Primera'
$ws.Range("J5").Value = 20
$ws.Range("K5").Value = 30000
$ws.Range("L5").Value = 30000
$ws.Range("M5").Value = 30000
$ws.Range("N5").Value = '$/caja 25 kilos'
$ws.Range("P5").Value = 1200
$ws.Range("Q5").Value = 25
$ws.Range("D6").Value = 44326
$ws.Range("H6").Value = 'Americana (o)'
$ws.Range("I6").Value = 'Primera'
$ws.Range("J6").Value = 15
$ws.Range("K6").Value = 30000
$ws.Range("L6").Value = 30000
$ws.Range("M6").Value = 30000
$ws.Range("N6").Value = '$/caja 25 kilos'
$ws.Range("P6").Value = 1200
$ws.Range("Q6").Value = 25
$ws.Range("D7").Value = 44460
$ws.Range("H7").Value = 'Americana (o)'
$ws.Range("I7").Value = 'Primera'
$ws.Range("J7").Value = 30
$ws.Range("K7").Value = 95000
$ws.Range("L7").Value = 95000
$ws.Range("M7").Value = 95000
$ws.Range("N7").Value = '$/caja 25 kilos'
$ws.Range("P7").Value = 3800
$ws.Range("Q7").Value = 25
$ws.Range("D8").Value = 44553
$ws.Range("H8").Value = 'Inferno'
$ws.Range("I8").Value = 'Primera'
$ws.Range("J8").Value = 35
$ws.Range("K8").Value = 45000
$ws.Range("L8").Value = 45000
$ws.Range("M8").Value = 45000
$ws.Range("N8").Value = '$/caja 25 kilos'
$ws.Range("P8").Value = 1800
$ws.Range("Q8").Value = 25
$ws.Range("D9").Value = 44425
$ws.Range("H9").Value = 'Americana (o)'
$ws.Range("I9").Value = 'Primera'
$ws.Range("J9").Value = 15
$ws.Range("K9").Value = 75000
$ws.Range("L9").Value = 75000
$ws.Range("M9").Value = 75000
$ws.Range("N9").Value = '$/caja 25 kilos'
$ws.Range("P9").Value = 3000
$ws.Range("Q9").Value = 25
$ws.Range("D11").Value = 44474
$ws.Range("H11").Value = 'Americana (o)'
$ws.Range("I11").Value = 'Primera'
$ws.Range("J11").Value = 18
$ws.Range("K11").Value = 100000
$ws.Range("L11").Value = 100000
$ws.Range("M11").Value = 100000
$ws.Range("N11").Value = '$/caja 25 kilos'
$ws.Range("P11").Value = 4000
$ws.Range("Q11").Value = 25
$ws.Range("D12").Value = 44446
$ws.Range("H12").Value = 'Americana (o)'
$ws.Range("I12").Value = 'Primera'
$ws.Range("J12").Value = 5
$ws.Range("K12").Value = 78000
$ws.Range("L12").Value = 78000
$ws.Range("M12").Value = 78000
$ws.Range("N12").Value = '$/caja 25 kilos'
$ws.Range("P12").Value = 3120
$ws.Range("Q12").Value = 25
$ws.Range("D13").Value = 44446
$ws.Range("H13").Value = 'Inferno'
$ws.Range("I13").Value = 'Primera'
$ws.Range("J13").Value = 4
$ws.Range("K13").Value = 80000
$ws.Range("L13").Value = 80000
$ws.Range("M13").Value = 80000
$ws.Range("N13").Value = '$/caja 15 kilos'
$ws.Range("P13").Value = 5333
$ws.Range("Q13").Value = 15
$ws.Range("D14").Value = 44340
$ws.Range("H14").Value = 'Americana (o)'
$ws.Range("I14").Value = 'Primera'
$ws.Range("J14").Value = 15
$ws.Range("K14").Value = 35000
$ws.Range("L14").Value = 35000
$ws.Range("M14").Value = 35000
$ws.Range("N14").Value = '$/caja 25 kilos'
$ws.Range("P14").Value = 1400
$ws.Range("Q14").Value = 25
$ws.Range("D15").Value = 44449
$ws.Range("H15").Value = 'Americana (o)'
$ws.Range("I15").Value = 'Primera'
$ws.Range("J15").Value = 25
$ws.Range("K15").Value = 80000
$ws.Range("L15").Value = 80000
$ws.Range("M15").Value = 80000
$ws.Range("N15").Value = '$/caja 25 kilos'
$ws.Range("P15").Value = 3200
$ws.Range("Q15").Value = 25
$ws.Range("D16").Value = 44449
$ws.Range("H16").Value = 'Americana (o)'
$ws.Range("I16").Value = 'Segunda'
$ws.Range("J16").Value = 20
$ws.Range("K16").Value = 75000
$ws.Range("L16").Value = 75000
$ws.Range("M16").Value = 75000
$ws.Range("N16").Value = '$/caja 15 kilos'
$ws.Range("P16").Value = 5000
$ws.Range("Q16").Value = 15
$ws.Range("D17").Value = 44193
$ws.Range("H17").Value = 'Americana (o)'
$ws.Range("I17").Value = 'Primera'
$ws.Range("J17").Value = 15
$ws.Range("K17").Value = 46000
$ws.Range("L17").Value = 46000
$ws.Range("M17").Value = 46000
$ws.Range("N17").Value = '$/caja 15 kilos'
$ws.Range("P17").Value = 3067
$ws.Range("Q17").Value = 15
$ws.Range("D18").Value = 44544
$ws.Range("H18").Value = 'Inferno'
$ws.Range("I18").Value = 'Primera'
$ws.Range("J18").Value = 12
$ws.Range("K18").Value = 35000
$ws.Range("L18").Value = 35000
$ws.Range("M18").Value = 35000
$ws.Range("N18").Value = '$/caja 25 kilos'
$ws.Range("P18").Value = 1400
$ws.Range("Q18").Value = 25
